$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "phoneNumber" test-case row (old row 64). This shifts
# all subsequent rows (old 65-81) up by one, carrying their values, styles
# and row heights along with them, which matches the target workbook where
# the sheet now ends at row 80 instead of row 81.
$ws.Rows(64).Delete()

# Fix a stray/incorrect "Testcase Objective" value that had been copy-pasted
# from a different test (it should match the actual steps recorded in the
# Steps column for this row, now row 73 after the deletion above).
$ws.Range("E73").Value = "wait(3);`nvalidate1;`nlink_Click(system_test_link);`nvalidate2;`nSelectTestToRun(VT300_081_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nvalidate4;`nvalidate5;"

# Restore the saved cursor/selection position on the sheet.
$ws.Range("D2").Select()
